# Fitur Multi Size & Gambar Detail
# Rename three header/column labels on the "barang" sheet:
#   id_category -> id_kategori
#   code_barang -> kode_barang
#   price       -> harga
# (the underlying data for every row is unchanged; only the header text
#  for columns B, C and F is renamed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "id_kategori"
$ws.Range("C1").Value = "kode_barang"
$ws.Range("F1").Value = "harga"

# Update the view state to match: selection moves from H7 to K1, and the
# window is scrolled so column G is the left-most visible column.
$ws.Activate()
$null = $ws.Range("K1").Select()

$win = $excel.ActiveWindow
try { $win.Left = 1750 } catch {}
try { $win.Top = 410 } catch {}
